$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "test"
